$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 269.54285
$ws.Range("I9").Value = 284.2069
$ws.Range("K9").Value = 284.2069
$ws.Range("M9").Value = -115.2069
$ws.Range("H11").Value = 53.375
$ws.Range("I11").Value = 53.375
$ws.Range("K11").Value = 53.375
$ws.Range("M11").Value = 86.625
$ws.Range("H19").Value = 705.9167
$ws.Range("I19").Value = 857.2857
$ws.Range("J19").Value = 494
$ws.Range("K19").Value = 857.2857
$ws.Range("L19").Value = 494
$ws.Range("M19").Value = -682.2857
$ws.Range("N19").Value = -844
$ws.Range("H33").Value = 411.7353
$ws.Range("I33").Value = 430.21875
$ws.Range("J33").Value = 116
$ws.Range("K33").Value = 430.21875
$ws.Range("L33").Value = 116
$ws.Range("M33").Value = -201.21875
$ws.Range("N33").Value = -574
$ws.Range("H80").Value = 949.7838
$ws.Range("I80").Value = 444.94116
$ws.Range("J80").Value = 1378.9
$ws.Range("K80").Value = 1334.82348
$ws.Range("L80").Value = 4136.700000000001
$ws.Range("M80").Value = -336.82348
$ws.Range("N80").Value = -6132.700000000001
$ws.Range("H83").Value = 949.7838
$ws.Range("I83").Value = 444.94116
$ws.Range("J83").Value = 1378.9
$ws.Range("K83").Value = 4004.47044
$ws.Range("L83").Value = 12410.1
$ws.Range("M83").Value = 987.5295599999999
$ws.Range("N83").Value = -22394.1
$ws.Range("H88").Value = 1993.5676
$ws.Range("I88").Value = 4375
$ws.Range("J88").Value = 1111.5555
$ws.Range("K88").Value = 4375
$ws.Range("L88").Value = 1111.5555
$ws.Range("M88").Value = -3969
$ws.Range("N88").Value = -1923.5555
$ws.Range("H91").Value = 1993.5676
$ws.Range("I91").Value = 4375
$ws.Range("J91").Value = 1111.5555
$ws.Range("K91").Value = 4375
$ws.Range("L91").Value = 1111.5555
$ws.Range("M91").Value = -2971
$ws.Range("N91").Value = -3919.5555
$ws.Range("H106").Value = 762.8570999999999
$ws.Range("I106").Value = 762.8570999999999
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 762.8570999999999
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H112").Value = 1117063.5
$ws.Range("I112").Value = 375
$ws.Range("J112").Value = 1563738.9
$ws.Range("K112").Value = 1125
$ws.Range("L112").Value = 4691216.699999999
$ws.Range("M112").Value = -17
$ws.Range("N112").Value = -4693432.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2419.8125
$ws.Range("I61").Value = 2210.5
$ws.Range("J61").Value = 3885
$ws.Range("K61").Value = 2210.5
$ws.Range("L61").Value = 3885
$ws.Range("M61").Value = -1998.5
$ws.Range("N61").Value = -4309
$ws.Range("H63").Value = 2278.125
$ws.Range("I63").Value = 2103.5715
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 2103.5715
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -1417.5715
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 2278.125
$ws.Range("I66").Value = 2103.5715
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 10517.8575
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -7085.8575
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 629202.2
$ws.Range("I74").Value = 5302.9165
$ws.Range("J74").Value = 2500900
$ws.Range("K74").Value = 5302.9165
$ws.Range("L74").Value = 2500900
$ws.Range("M74").Value = -4428.9165
$ws.Range("N74").Value = -2502648
$ws.Range("H77").Value = 629202.2
$ws.Range("I77").Value = 5302.9165
$ws.Range("J77").Value = 2500900
$ws.Range("K77").Value = 26514.5825
$ws.Range("L77").Value = 12504500
$ws.Range("M77").Value = -22146.5825
$ws.Range("N77").Value = -12513236
$ws.Range("H102").Value = 1546.5555
$ws.Range("I102").Value = 1199
$ws.Range("K102").Value = 1199
$ws.Range("M102").Value = 423
$ws.Range("H134").Value = 48571.43
$ws.Range("J134").Value = 48571.43
$ws.Range("L134").Value = 48571.43
$ws.Range("N134").Value = -58711.43
$ws.Range("H136").Value = 2419.8125
$ws.Range("I136").Value = 2210.5
$ws.Range("J136").Value = 3885
$ws.Range("K136").Value = 6631.5
$ws.Range("L136").Value = 11655
$ws.Range("M136").Value = -4081.5
$ws.Range("N136").Value = -16755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2458.5454
$ws.Range("I20").Value = 1272.8096
$ws.Range("J20").Value = 4533.5835
$ws.Range("K20").Value = 1272.8096
$ws.Range("L20").Value = 4533.5835
$ws.Range("M20").Value = -1025.8096
$ws.Range("N20").Value = -5027.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1700.5
$ws.Range("I35").Value = 1700.5
$ws.Range("K35").Value = 1700.5
$ws.Range("M35").Value = -1406.5
$ws.Range("H134").Value = 2702.85
$ws.Range("I134").Value = 2814.4
$ws.Range("J134").Value = 1922
$ws.Range("K134").Value = 8443.200000000001
$ws.Range("L134").Value = 5766
$ws.Range("M134").Value = -5908.200000000001
$ws.Range("N134").Value = -10836

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 123770.125
$ws.Range("I2").Value = 198017.2
$ws.Range("K2").Value = 1188103.2
$ws.Range("M2").Value = -1187990.2
$ws.Range("H5").Value = 678.2154
$ws.Range("I5").Value = 514.375
$ws.Range("K5").Value = 1543.125
$ws.Range("M5").Value = -1431.125
$ws.Range("H68").Value = 1597.125
$ws.Range("J68").Value = 1746.1666
$ws.Range("L68").Value = 5238.4998
$ws.Range("N68").Value = -6860.4998
$ws.Range("H71").Value = 1597.125
$ws.Range("J71").Value = 1746.1666
$ws.Range("L71").Value = 15715.4994
$ws.Range("N71").Value = -23827.4994
$ws.Range("H131").Value = 839.5294
$ws.Range("I131").Value = 514.7895
$ws.Range("J131").Value = 1250.8667
$ws.Range("K131").Value = 1544.3685
$ws.Range("L131").Value = 3752.6001
$ws.Range("M131").Value = 3495.6315
$ws.Range("N131").Value = -13832.6001
$ws.Range("H132").Value = 1009.25
$ws.Range("I132").Value = 501.73685
$ws.Range("J132").Value = 1576.4706
$ws.Range("K132").Value = 4515.63165
$ws.Range("L132").Value = 14188.2354
$ws.Range("M132").Value = -1985.63165
$ws.Range("N132").Value = -19248.2354
$ws.Range("H135").Value = 678.2154
$ws.Range("I135").Value = 514.375
$ws.Range("K135").Value = 4629.375
$ws.Range("M135").Value = -2094.375
$ws.Range("H141").Value = 1432.5294
$ws.Range("I141").Value = 668.7368
$ws.Range("J141").Value = 2400
$ws.Range("K141").Value = 2006.2104
$ws.Range("L141").Value = 7200
$ws.Range("M141").Value = 3173.7896
$ws.Range("N141").Value = -17560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 217.46666
$ws.Range("I2").Value = 151.66667
$ws.Range("J2").Value = 316.16666
$ws.Range("K2").Value = 151.66667
$ws.Range("L2").Value = 316.16666
$ws.Range("M2").Value = -38.66667000000001
$ws.Range("N2").Value = -542.16666
$ws.Range("H113").Value = 1732.3043
$ws.Range("I113").Value = 1060.0667
$ws.Range("K113").Value = 1060.0667
$ws.Range("M113").Value = 1109.9333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3015
$ws.Range("J32").Value = 3015
$ws.Range("L32").Value = 3015
$ws.Range("N32").Value = -3649
$ws.Range("H55").Value = 280.7647
$ws.Range("I55").Value = 289.57895
$ws.Range("J55").Value = 269.6
$ws.Range("K55").Value = 289.57895
$ws.Range("L55").Value = 269.6
$ws.Range("M55").Value = -116.57895
$ws.Range("N55").Value = -615.6
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2960.842
$ws.Range("I136").Value = 996.6829
$ws.Range("J136").Value = 7994
$ws.Range("K136").Value = 2990.0487
$ws.Range("L136").Value = 23982
$ws.Range("M136").Value = -440.0487000000003
$ws.Range("N136").Value = -29082
